$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 14:57:55"
$ws1.Cells.Item(3,1).Value = "Total filas: 314"
$ws1.Cells.Item(14,1).Value = "04:44:55"
$ws1.Cells.Item(14,3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(14,4).Value = 2
$ws1.Cells.Item(15,1).Value = "03:46:12"
$ws1.Cells.Item(15,3).Value = "215A_EL PATO"
$ws1.Cells.Item(15,4).Value = 60
$ws1.Cells.Item(43,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(44,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(80,1).Value = "07:31:04"
$ws1.Cells.Item(80,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(80,4).Value = 28
$ws1.Cells.Item(81,1).Value = "06:28:32"
$ws1.Cells.Item(81,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(81,4).Value = 91
$ws1.Cells.Item(151,1).Value = "08:41:14"
$ws1.Cells.Item(151,3).Value = "10_OLMOS"
$ws1.Cells.Item(151,4).Value = 113
$ws1.Cells.Item(152,1).Value = "09:57:03"
$ws1.Cells.Item(152,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(152,4).Value = 37
$ws1.Cells.Item(200,1).Value = "10:51:31"
$ws1.Cells.Item(200,3).Value = "215C_EL PATO"
$ws1.Cells.Item(200,4).Value = 98
$ws1.Cells.Item(201,1).Value = "12:29:23"
$ws1.Cells.Item(201,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(201,4).Value = 0
$ws1.Cells.Item(228,1).Value = "11:45:06"
$ws1.Cells.Item(228,3).Value = "215_ALUAR"
$ws1.Cells.Item(228,4).Value = 86
$ws1.Cells.Item(229,1).Value = "12:29:23"
$ws1.Cells.Item(229,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(229,4).Value = 42
$ws1.Cells.Item(236,3).Value = "215A_EL PATO"
$ws1.Cells.Item(237,3).Value = "10_OLMOS"
$ws1.Cells.Item(249,1).Value = "12:29:23"
$ws1.Cells.Item(249,3).Value = "215C_LA PLATA"
$ws1.Cells.Item(249,4).Value = 85
$ws1.Cells.Item(250,1).Value = "11:58:34"
$ws1.Cells.Item(250,3).Value = "225_GOMEZ"
$ws1.Cells.Item(250,4).Value = 116
$ws1.Cells.Item(281,1).Value = "14:57:55"
$ws1.Cells.Item(281,2).Value = "14:58"
$ws1.Cells.Item(281,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(281,4).Value = 1
$ws1.Cells.Item(282,1).Value = "14:43:48"
$ws1.Cells.Item(282,2).Value = "15:00"
$ws1.Cells.Item(282,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(282,4).Value = 17
$ws1.Cells.Item(283,1).Value = "13:35:25"
$ws1.Cells.Item(283,2).Value = "15:01"
$ws1.Cells.Item(283,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(283,4).Value = 86
$ws1.Cells.Item(284,1).Value = "13:54:15"
$ws1.Cells.Item(284,2).Value = "15:02"
$ws1.Cells.Item(284,3).Value = "215A_LA PLATA"
$ws1.Cells.Item(284,4).Value = 68
$ws1.Cells.Item(285,1).Value = "14:57:55"
$ws1.Cells.Item(285,2).Value = "15:02"
$ws1.Cells.Item(285,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(285,4).Value = 5
$ws1.Cells.Item(286,2).Value = "15:03"
$ws1.Cells.Item(286,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(286,4).Value = 40
$ws1.Cells.Item(287,1).Value = "13:35:25"
$ws1.Cells.Item(287,3).Value = "14_ABASTO"
$ws1.Cells.Item(287,4).Value = 89
$ws1.Cells.Item(288,1).Value = "14:23:38"
$ws1.Cells.Item(288,2).Value = "15:04"
$ws1.Cells.Item(288,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(288,4).Value = 41
$ws1.Cells.Item(289,1).Value = "14:43:48"
$ws1.Cells.Item(289,2).Value = "15:04"
$ws1.Cells.Item(289,3).Value = "215A_LA PLATA"
$ws1.Cells.Item(289,4).Value = 21
$ws1.Cells.Item(290,1).Value = "13:54:15"
$ws1.Cells.Item(290,2).Value = "15:05"
$ws1.Cells.Item(290,3).Value = "14_ABASTO"
$ws1.Cells.Item(290,4).Value = 71
$ws1.Cells.Item(291,1).Value = "13:35:25"
$ws1.Cells.Item(291,2).Value = "15:17"
$ws1.Cells.Item(291,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(291,4).Value = 102
$ws1.Cells.Item(292,1).Value = "13:35:25"
$ws1.Cells.Item(292,2).Value = "15:24"
$ws1.Cells.Item(292,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(292,4).Value = 109
$ws1.Cells.Item(293,1).Value = "14:23:38"
$ws1.Cells.Item(293,2).Value = "15:24"
$ws1.Cells.Item(293,4).Value = 61
$ws1.Cells.Item(294,1).Value = "13:35:25"
$ws1.Cells.Item(294,3).Value = "215C_EL PATO"
$ws1.Cells.Item(294,4).Value = 110
$ws1.Cells.Item(295,2).Value = "15:25"
$ws1.Cells.Item(295,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(295,4).Value = 91
$ws1.Cells.Item(296,1).Value = "13:54:15"
$ws1.Cells.Item(296,2).Value = "15:25"
$ws1.Cells.Item(296,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(296,4).Value = 91
$ws1.Cells.Item(297,1).Value = "14:57:55"
$ws1.Cells.Item(297,2).Value = "15:33"
$ws1.Cells.Item(297,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(297,4).Value = 36
$ws1.Cells.Item(298,1).Value = "13:54:15"
$ws1.Cells.Item(298,2).Value = "15:36"
$ws1.Cells.Item(298,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(298,4).Value = 102
$ws1.Cells.Item(299,2).Value = "15:44"
$ws1.Cells.Item(299,3).Value = "14_ABASTO"
$ws1.Cells.Item(299,4).Value = 81
$ws1.Cells.Item(300,1).Value = "14:23:38"
$ws1.Cells.Item(300,2).Value = "15:45"
$ws1.Cells.Item(300,3).Value = "215C_LA PLATA"
$ws1.Cells.Item(300,4).Value = 82
$ws1.Cells.Item(301,1).Value = "14:57:55"
$ws1.Cells.Item(301,2).Value = "15:51"
$ws1.Cells.Item(301,3).Value = "215C_LA PLATA"
$ws1.Cells.Item(301,4).Value = 54
$ws1.Cells.Item(302,1).Value = "14:43:48"
$ws1.Cells.Item(302,2).Value = "15:53"
$ws1.Cells.Item(302,3).Value = "215C_LA PLATA"
$ws1.Cells.Item(302,4).Value = 70
$ws1.Cells.Item(303,1).Value = "14:23:38"
$ws1.Cells.Item(303,2).Value = "15:55"
$ws1.Cells.Item(303,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(303,4).Value = 92
$ws1.Cells.Item(304,1).Value = "14:23:38"
$ws1.Cells.Item(304,2).Value = "15:56"
$ws1.Cells.Item(304,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(304,4).Value = 73
$ws1.Cells.Item(305,1).Value = "14:23:38"
$ws1.Cells.Item(305,2).Value = "16:01"
$ws1.Cells.Item(305,3).Value = "15_ABASTO"
$ws1.Cells.Item(305,4).Value = 98
$ws1.Cells.Item(306,1).Value = "14:57:55"
$ws1.Cells.Item(306,2).Value = "16:16"
$ws1.Cells.Item(306,3).Value = "10_OLMOS"
$ws1.Cells.Item(306,4).Value = 79
$ws1.Cells.Item(307,1).Value = "14:23:38"
$ws1.Cells.Item(307,2).Value = "16:20"
$ws1.Cells.Item(307,3).Value = "10_OLMOS"
$ws1.Cells.Item(307,4).Value = 117
$ws1.Cells.Item(308,1).Value = "14:57:55"
$ws1.Cells.Item(308,2).Value = "16:22"
$ws1.Cells.Item(308,3).Value = "15_ABASTO"
$ws1.Cells.Item(308,4).Value = 85
$ws1.Cells.Item(309,2).Value = "16:24"
$ws1.Cells.Item(309,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(309,4).Value = 101
$ws1.Cells.Item(310,1).Value = "14:43:48"
$ws1.Cells.Item(310,2).Value = "16:24"
$ws1.Cells.Item(310,3).Value = "215_ALUAR"
$ws1.Cells.Item(310,4).Value = 101
$ws1.Cells.Item(310,5).Value = "LP1912"
$ws1.Cells.Item(311,1).Value = "14:43:48"
$ws1.Cells.Item(311,2).Value = "16:25"
$ws1.Cells.Item(311,3).Value = "215B_LP-P MOR-1 Y 57"
$ws1.Cells.Item(311,4).Value = 102
$ws1.Cells.Item(311,5).Value = "LP1912"
$ws1.Cells.Item(312,1).Value = "14:43:48"
$ws1.Cells.Item(312,2).Value = "16:31"
$ws1.Cells.Item(312,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(312,4).Value = 108
$ws1.Cells.Item(312,5).Value = "LP1912"
$ws1.Cells.Item(313,1).Value = "14:43:48"
$ws1.Cells.Item(313,2).Value = "16:33"
$ws1.Cells.Item(313,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(313,4).Value = 110
$ws1.Cells.Item(313,5).Value = "LP1912"
$ws1.Cells.Item(314,1).Value = "14:57:55"
$ws1.Cells.Item(314,2).Value = "16:33"
$ws1.Cells.Item(314,3).Value = "225_GOMEZ"
$ws1.Cells.Item(314,4).Value = 96
$ws1.Cells.Item(314,5).Value = "LP1912"
$ws1.Cells.Item(315,1).Value = "14:43:48"
$ws1.Cells.Item(315,2).Value = "16:34"
$ws1.Cells.Item(315,3).Value = "225_GOMEZ"
$ws1.Cells.Item(315,4).Value = 111
$ws1.Cells.Item(315,5).Value = "LP1912"
$ws1.Cells.Item(316,1).Value = "14:43:48"
$ws1.Cells.Item(316,2).Value = "16:36"
$ws1.Cells.Item(316,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(316,4).Value = 113
$ws1.Cells.Item(316,5).Value = "LP1912"
$ws1.Cells.Item(317,1).Value = "14:57:55"
$ws1.Cells.Item(317,2).Value = "16:41"
$ws1.Cells.Item(317,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(317,4).Value = 104
$ws1.Cells.Item(317,5).Value = "LP1912"
$ws1.Cells.Item(318,1).Value = "14:57:55"
$ws1.Cells.Item(318,2).Value = "16:53"
$ws1.Cells.Item(318,3).Value = "10_OLMOS"
$ws1.Cells.Item(318,4).Value = 116
$ws1.Cells.Item(318,5).Value = "LP1912"
$ws1.Cells.Item(319,1).Value = "14:57:55"
$ws1.Cells.Item(319,2).Value = "16:56"
$ws1.Cells.Item(319,3).Value = "215C_EL PATO"
$ws1.Cells.Item(319,4).Value = 119
$ws1.Cells.Item(319,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 14:57:55"
$ws2.Cells.Item(3,1).Value = "Total filas: 82"
$ws2.Cells.Item(83,1).Value = "14:57:55"
$ws2.Cells.Item(83,2).Value = "15:51"
$ws2.Cells.Item(83,4).Value = 54
$ws2.Cells.Item(84,2).Value = "15:53"
$ws2.Cells.Item(84,3).Value = "215C_LA PLATA"
$ws2.Cells.Item(84,4).Value = 70
$ws2.Cells.Item(85,2).Value = "16:24"
$ws2.Cells.Item(85,3).Value = "215_ALUAR"
$ws2.Cells.Item(85,4).Value = 101
$ws2.Cells.Item(86,1).Value = "14:43:48"
$ws2.Cells.Item(86,2).Value = "16:25"
$ws2.Cells.Item(86,3).Value = "215B_LP-P MOR-1 Y 57"
$ws2.Cells.Item(86,4).Value = 102
$ws2.Cells.Item(86,5).Value = "LP1912"
$ws2.Cells.Item(87,1).Value = "14:57:55"
$ws2.Cells.Item(87,2).Value = "16:56"
$ws2.Cells.Item(87,3).Value = "215C_EL PATO"
$ws2.Cells.Item(87,4).Value = 119
$ws2.Cells.Item(87,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 14:57:55"
$ws3.Cells.Item(3,1).Value = "Total filas: 314"
$ws3.Cells.Item(14,1).Value = "04:44:55"
$ws3.Cells.Item(14,3).Value = "215_EL PELIGRO"
$ws3.Cells.Item(14,4).Value = 2
$ws3.Cells.Item(15,1).Value = "03:46:12"
$ws3.Cells.Item(15,3).Value = "215A_EL PATO"
$ws3.Cells.Item(15,4).Value = 60
$ws3.Cells.Item(43,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(44,3).Value = "17X38_ROMERO"
$ws3.Cells.Item(80,1).Value = "07:31:04"
$ws3.Cells.Item(80,3).Value = "23_HERNANDEZ"
$ws3.Cells.Item(80,4).Value = 28
$ws3.Cells.Item(81,1).Value = "06:28:32"
$ws3.Cells.Item(81,3).Value = "11_ETCHEVERRY"
$ws3.Cells.Item(81,4).Value = 91
$ws3.Cells.Item(151,1).Value = "08:41:14"
$ws3.Cells.Item(151,3).Value = "10_OLMOS"
$ws3.Cells.Item(151,4).Value = 113
$ws3.Cells.Item(152,1).Value = "09:57:03"
$ws3.Cells.Item(152,3).Value = "23_HERNANDEZ"
$ws3.Cells.Item(152,4).Value = 37
$ws3.Cells.Item(200,1).Value = "10:51:31"
$ws3.Cells.Item(200,3).Value = "215C_EL PATO"
$ws3.Cells.Item(200,4).Value = 98
$ws3.Cells.Item(201,1).Value = "12:29:23"
$ws3.Cells.Item(201,3).Value = "11_ETCHEVERRY"
$ws3.Cells.Item(201,4).Value = 0
$ws3.Cells.Item(228,1).Value = "11:45:06"
$ws3.Cells.Item(228,3).Value = "215_ALUAR"
$ws3.Cells.Item(228,4).Value = 86
$ws3.Cells.Item(229,1).Value = "12:29:23"
$ws3.Cells.Item(229,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(229,4).Value = 42
$ws3.Cells.Item(236,3).Value = "215A_EL PATO"
$ws3.Cells.Item(237,3).Value = "10_OLMOS"
$ws3.Cells.Item(249,1).Value = "12:29:23"
$ws3.Cells.Item(249,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(249,4).Value = 85
$ws3.Cells.Item(250,1).Value = "11:58:34"
$ws3.Cells.Item(250,3).Value = "225_GOMEZ"
$ws3.Cells.Item(250,4).Value = 116
$ws3.Cells.Item(281,1).Value = "14:57:55"
$ws3.Cells.Item(281,2).Value = "14:58"
$ws3.Cells.Item(281,3).Value = "16_P MOR-SANTA ANA"
$ws3.Cells.Item(281,4).Value = 1
$ws3.Cells.Item(282,1).Value = "14:43:48"
$ws3.Cells.Item(282,2).Value = "15:00"
$ws3.Cells.Item(282,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(282,4).Value = 17
$ws3.Cells.Item(283,1).Value = "13:35:25"
$ws3.Cells.Item(283,2).Value = "15:01"
$ws3.Cells.Item(283,3).Value = "81_EL PELIGRO"
$ws3.Cells.Item(283,4).Value = 86
$ws3.Cells.Item(284,1).Value = "13:54:15"
$ws3.Cells.Item(284,2).Value = "15:02"
$ws3.Cells.Item(284,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(284,4).Value = 68
$ws3.Cells.Item(285,1).Value = "14:57:55"
$ws3.Cells.Item(285,2).Value = "15:02"
$ws3.Cells.Item(285,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(285,4).Value = 5
$ws3.Cells.Item(286,2).Value = "15:03"
$ws3.Cells.Item(286,3).Value = "23_HERNANDEZ"
$ws3.Cells.Item(286,4).Value = 40
$ws3.Cells.Item(287,1).Value = "13:35:25"
$ws3.Cells.Item(287,3).Value = "14_ABASTO"
$ws3.Cells.Item(287,4).Value = 89
$ws3.Cells.Item(288,1).Value = "14:23:38"
$ws3.Cells.Item(288,2).Value = "15:04"
$ws3.Cells.Item(288,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(288,4).Value = 41
$ws3.Cells.Item(289,1).Value = "14:43:48"
$ws3.Cells.Item(289,2).Value = "15:04"
$ws3.Cells.Item(289,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(289,4).Value = 21
$ws3.Cells.Item(290,1).Value = "13:54:15"
$ws3.Cells.Item(290,2).Value = "15:05"
$ws3.Cells.Item(290,3).Value = "14_ABASTO"
$ws3.Cells.Item(290,4).Value = 71
$ws3.Cells.Item(291,1).Value = "13:35:25"
$ws3.Cells.Item(291,2).Value = "15:17"
$ws3.Cells.Item(291,3).Value = "16_P MOR-SANTA ANA"
$ws3.Cells.Item(291,4).Value = 102
$ws3.Cells.Item(292,1).Value = "13:35:25"
$ws3.Cells.Item(292,2).Value = "15:24"
$ws3.Cells.Item(292,3).Value = "11_ETCHEVERRY"
$ws3.Cells.Item(292,4).Value = 109
$ws3.Cells.Item(293,1).Value = "14:23:38"
$ws3.Cells.Item(293,2).Value = "15:24"
$ws3.Cells.Item(293,4).Value = 61
$ws3.Cells.Item(294,1).Value = "13:35:25"
$ws3.Cells.Item(294,3).Value = "215C_EL PATO"
$ws3.Cells.Item(294,4).Value = 110
$ws3.Cells.Item(295,2).Value = "15:25"
$ws3.Cells.Item(295,3).Value = "11_ETCHEVERRY"
$ws3.Cells.Item(295,4).Value = 91
$ws3.Cells.Item(296,1).Value = "13:54:15"
$ws3.Cells.Item(296,2).Value = "15:25"
$ws3.Cells.Item(296,3).Value = "16_P MOR-SANTA ANA"
$ws3.Cells.Item(296,4).Value = 91
$ws3.Cells.Item(297,1).Value = "14:57:55"
$ws3.Cells.Item(297,2).Value = "15:33"
$ws3.Cells.Item(297,3).Value = "16_SANTA ANA"
$ws3.Cells.Item(297,4).Value = 36
$ws3.Cells.Item(298,1).Value = "13:54:15"
$ws3.Cells.Item(298,2).Value = "15:36"
$ws3.Cells.Item(298,3).Value = "17X38_ROMERO"
$ws3.Cells.Item(298,4).Value = 102
$ws3.Cells.Item(299,2).Value = "15:44"
$ws3.Cells.Item(299,3).Value = "14_ABASTO"
$ws3.Cells.Item(299,4).Value = 81
$ws3.Cells.Item(300,1).Value = "14:23:38"
$ws3.Cells.Item(300,2).Value = "15:45"
$ws3.Cells.Item(300,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(300,4).Value = 82
$ws3.Cells.Item(301,1).Value = "14:57:55"
$ws3.Cells.Item(301,2).Value = "15:51"
$ws3.Cells.Item(301,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(301,4).Value = 54
$ws3.Cells.Item(302,1).Value = "14:43:48"
$ws3.Cells.Item(302,2).Value = "15:53"
$ws3.Cells.Item(302,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(302,4).Value = 70
$ws3.Cells.Item(303,1).Value = "14:23:38"
$ws3.Cells.Item(303,2).Value = "15:55"
$ws3.Cells.Item(303,3).Value = "27_EL RETIRO"
$ws3.Cells.Item(303,4).Value = 92
$ws3.Cells.Item(304,1).Value = "14:23:38"
$ws3.Cells.Item(304,2).Value = "15:56"
$ws3.Cells.Item(304,3).Value = "27_EL RETIRO"
$ws3.Cells.Item(304,4).Value = 73
$ws3.Cells.Item(305,1).Value = "14:23:38"
$ws3.Cells.Item(305,2).Value = "16:01"
$ws3.Cells.Item(305,3).Value = "15_ABASTO"
$ws3.Cells.Item(305,4).Value = 98
$ws3.Cells.Item(306,1).Value = "14:57:55"
$ws3.Cells.Item(306,2).Value = "16:16"
$ws3.Cells.Item(306,3).Value = "10_OLMOS"
$ws3.Cells.Item(306,4).Value = 79
$ws3.Cells.Item(307,1).Value = "14:23:38"
$ws3.Cells.Item(307,2).Value = "16:20"
$ws3.Cells.Item(307,3).Value = "10_OLMOS"
$ws3.Cells.Item(307,4).Value = 117
$ws3.Cells.Item(308,1).Value = "14:57:55"
$ws3.Cells.Item(308,2).Value = "16:22"
$ws3.Cells.Item(308,3).Value = "15_ABASTO"
$ws3.Cells.Item(308,4).Value = 85
$ws3.Cells.Item(309,2).Value = "16:24"
$ws3.Cells.Item(309,3).Value = "11_ETCHEVERRY"
$ws3.Cells.Item(309,4).Value = 101
$ws3.Cells.Item(310,1).Value = "14:43:48"
$ws3.Cells.Item(310,2).Value = "16:24"
$ws3.Cells.Item(310,3).Value = "215_ALUAR"
$ws3.Cells.Item(310,4).Value = 101
$ws3.Cells.Item(310,5).Value = "LP1912"
$ws3.Cells.Item(311,1).Value = "14:43:48"
$ws3.Cells.Item(311,2).Value = "16:25"
$ws3.Cells.Item(311,3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(311,4).Value = 102
$ws3.Cells.Item(311,5).Value = "LP1912"
$ws3.Cells.Item(312,1).Value = "14:43:48"
$ws3.Cells.Item(312,2).Value = "16:31"
$ws3.Cells.Item(312,3).Value = "16_P MOR-SANTA ANA"
$ws3.Cells.Item(312,4).Value = 108
$ws3.Cells.Item(312,5).Value = "LP1912"
$ws3.Cells.Item(313,1).Value = "14:43:48"
$ws3.Cells.Item(313,2).Value = "16:33"
$ws3.Cells.Item(313,3).Value = "23_HERNANDEZ"
$ws3.Cells.Item(313,4).Value = 110
$ws3.Cells.Item(313,5).Value = "LP1912"
$ws3.Cells.Item(314,1).Value = "14:57:55"
$ws3.Cells.Item(314,2).Value = "16:33"
$ws3.Cells.Item(314,3).Value = "225_GOMEZ"
$ws3.Cells.Item(314,4).Value = 96
$ws3.Cells.Item(314,5).Value = "LP1912"
$ws3.Cells.Item(315,1).Value = "14:43:48"
$ws3.Cells.Item(315,2).Value = "16:34"
$ws3.Cells.Item(315,3).Value = "225_GOMEZ"
$ws3.Cells.Item(315,4).Value = 111
$ws3.Cells.Item(315,5).Value = "LP1912"
$ws3.Cells.Item(316,1).Value = "14:43:48"
$ws3.Cells.Item(316,2).Value = "16:36"
$ws3.Cells.Item(316,3).Value = "17X38_ROMERO"
$ws3.Cells.Item(316,4).Value = 113
$ws3.Cells.Item(316,5).Value = "LP1912"
$ws3.Cells.Item(317,1).Value = "14:57:55"
$ws3.Cells.Item(317,2).Value = "16:41"
$ws3.Cells.Item(317,3).Value = "17X38_ROMERO"
$ws3.Cells.Item(317,4).Value = 104
$ws3.Cells.Item(317,5).Value = "LP1912"
$ws3.Cells.Item(318,1).Value = "14:57:55"
$ws3.Cells.Item(318,2).Value = "16:53"
$ws3.Cells.Item(318,3).Value = "10_OLMOS"
$ws3.Cells.Item(318,4).Value = 116
$ws3.Cells.Item(318,5).Value = "LP1912"
$ws3.Cells.Item(319,1).Value = "14:57:55"
$ws3.Cells.Item(319,2).Value = "16:56"
$ws3.Cells.Item(319,3).Value = "215C_EL PATO"
$ws3.Cells.Item(319,4).Value = 119
$ws3.Cells.Item(319,5).Value = "LP1912"
